# Applies the "outputs" -> "targets" wording tweak and relocates the
# auto-managed "_GoBack" bookmark to the end of that sentence (mirroring
# Word's own behaviour of stamping _GoBack at the last edited spot).

$d = $word.ActiveDocument

# --- 1. Swap "outputs" for "targets" in "... full-sized images for outputs ..." ---
$contentText = $d.Content.Text
$phraseIdx = $contentText.IndexOf("images for outputs")
$outputsIdx = $contentText.IndexOf("outputs", $phraseIdx)
$rOutputs = $d.Range($outputsIdx, $outputsIdx + "outputs".Length)
$rOutputs.Text = "targets"

# The text substitution above coalesces the touched paragraph's
# identically-formatted runs into one. Nudge formatting off/on (a
# no-op visually) on the runs that must stay distinct so the engine
# re-splits them back into separate <w:r> elements, matching how the
# surrounding, untouched runs ("full-sized" and the trailing
# " in both the training and validation sets.") looked originally.
$refreshedText = $d.Content.Text
$fullSizedIdx = $refreshedText.IndexOf("full-sized")
$rFullSized = $d.Range($fullSizedIdx, $fullSizedIdx + "full-sized".Length)
$rFullSized.Font.Bold = $true
$rFullSized.Font.Bold = $false

$targetsIdx = $refreshedText.IndexOf("targets", $fullSizedIdx)
$rTargets = $d.Range($targetsIdx, $targetsIdx + "targets".Length)
$rTargets.Font.Bold = $true
$rTargets.Font.Bold = $false

# --- 2. Move the "_GoBack" bookmark so it sits right after the sentence
#        that was just edited (". validation sets.") instead of after the
#        later "Mining / Learning from the data" heading. ---
$latestText = $d.Content.Text
$marker = "in both the training and validation sets."
$markerIdx = $latestText.IndexOf($marker)
$afterMarker = $markerIdx + $marker.Length
$rGoBack = $d.Range($afterMarker, $afterMarker)
$d.Bookmarks.Add("_GoBack", $rGoBack)
